$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: columns D1:T1 get reordered + new metric names ---
# (B1/C1 "series"/"spec" and A-column index are unchanged)
$headerVals = @(
    'rmse_thetass2_spme [milli]', 'rmse_phie2_spme [mV]', 'rmse_if2_spme [mA]', 'rmse_vcell_spme [mV]', 'rmse_thetass2_hybrid_without_delta [milli]', 'rmse_phie2_hybrid_without_delta [mV]', 'rmse_if2_hybrid_without_delta [mA]', 'rmse_vcell_hybrid_without_delta [mV]', 'rmse_thetass2_hybrid [milli]', 'rmse_phie2_hybrid [mV]', 'rmse_if2_hybrid [mA]', 'rmse_vcell_hybrid [mV]', 'rmse_thetass2_hybrid_mtau [milli]', 'rmse_phie2_hybrid_mtau [mV]', 'rmse_if2_hybrid_mtau [mA]', 'rmse_vcell_hybrid_mtau [mV]', 'rmse_vcell_rom [mV]'
)
$headerRow = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt $headerVals.Length; $i++) { $headerRow[0,$i] = $headerVals[$i] }
$ws.Range("D1:T1").Value = $headerRow

# New columns P:T need the same header style (bold, bordered, centered) as D:O
$ws.Range("D1").Copy()
$ws.Range("P1:T1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-8: columns D:S filled with recomputed RMSE values using tau from HRA ---
# Column T (rmse_vcell_rom) has no data, same as the old empty placeholder column.
$row2 = New-Object 'object[,]' 1,16
$row2Vals = @(43.7358703437796, 2.978995576824035, 56.67263062870067, 18.78470331533877, 2.293724908164353, 0.5341372674866498, 4.105678872680572, 1.251145118575022, 0.4007139611697978, 0.2897978153470505, 1.593423994716917, 0.4187197343067486, 0.3368233357377967, 0.196474083314284, 6.329179078424049, 0.3583155051906138)
for ($i = 0; $i -lt 16; $i++) { $row2[0,$i] = $row2Vals[$i] }
$ws.Range("D2:S2").Value = $row2

$row3 = New-Object 'object[,]' 1,16
$row3Vals = @(15.034764981421, 0.2532317906993753, 5.177685837604972, 5.003483852125205, 1.203935038555829, 0.4113714936147121, 20.62389045872475, 2.30295306187848, 0.2069777797965604, 0.1072690026844791, 0.5116107823734578, 0.1344227326742837, 0.2283399580241873, 0.125320613384622, 1.546829424739875, 0.06535990562676323)
for ($i = 0; $i -lt 16; $i++) { $row3[0,$i] = $row3Vals[$i] }
$ws.Range("D3:S3").Value = $row3

$row4 = New-Object 'object[,]' 1,16
$row4Vals = @(14.87308574524844, 0.4996402611546773, 9.120512990525478, 4.839307579931836, 3.163185622924866, 1.276481120222673, 5.427736819300264, 1.276188029148837, 0.2933831423953658, 0.07732002857186196, 1.36692736406588, 0.3322589499956418, 0.1338687111506674, 0.07512451265799253, 0.9090713512831607, 0.1371804276628667)
for ($i = 0; $i -lt 16; $i++) { $row4[0,$i] = $row4Vals[$i] }
$ws.Range("D4:S4").Value = $row4

$row5 = New-Object 'object[,]' 1,16
$row5Vals = @(46.71988130681122, 4.142925913006523, 77.16998448259233, 17.3959451475703, 3.569198843582928, 0.8915088172082858, 100.1547995140992, 7.384130577310017, 2.180694896785567, 0.3374207023560365, 8.577256505218603, 1.200392199277645, 2.779666701453919, 0.2102026984955428, 2.531300974628592, 4.126239001850039)
for ($i = 0; $i -lt 16; $i++) { $row5[0,$i] = $row5Vals[$i] }
$ws.Range("D5:S5").Value = $row5

$row6 = New-Object 'object[,]' 1,16
$row6Vals = @(18.02405058424596, 2.912686349415983, 52.79000399964171, 9.631803000573209, 13.24863715480514, 2.019695606016681, 40.54638496657045, 1.3146297405899, 0.578990546874148, 0.1535400392487568, 1.566693796553724, 0.2424488633090422, 0.8536955184065727, 0.7135600085984358, 3.642050367788539, 0.7065926621248783)
for ($i = 0; $i -lt 16; $i++) { $row6[0,$i] = $row6Vals[$i] }
$ws.Range("D6:S6").Value = $row6

$row7 = New-Object 'object[,]' 1,16
$row7Vals = @(18.20383718240869, 2.985287828596064, 54.3721708206039, 9.789901513567713, 11.76150623102837, 2.737432530400332, 43.35040114086348, 3.431444068153565, 0.642965281743617, 0.2531807654505578, 2.389374947162722, 0.3775667471293523, 1.46939076461661, 0.2291985881086432, 5.982585160341862, 1.111551481148627)
for ($i = 0; $i -lt 16; $i++) { $row7[0,$i] = $row7Vals[$i] }
$ws.Range("D7:S7").Value = $row7

$row8 = New-Object 'object[,]' 1,16
$row8Vals = @(23.53428404124674, 4.771704987708077, 102.6347181812276, 14.94918195657381, 14.07311921809941, 2.0945997301829, 40.66663299748471, 6.502474003138566, 1.184083850872454, 0.7110283462933842, 9.79870686147686, 1.529324178351341, 2.694588622602008, 0.8515649338313767, 9.122009511504032, 2.262576439384274)
for ($i = 0; $i -lt 16; $i++) { $row8[0,$i] = $row8Vals[$i] }
$ws.Range("D8:S8").Value = $row8

